$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "67.173.97"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.932.13"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.63%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "470.39"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +8.35%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.60"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.78%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.628"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.50%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.167"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +8.84%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0000342"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.68%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "43.42"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "10.47"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.555.65"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.57%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.17"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.15%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.924.61"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("E17").Value = "  +0.28%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "19.87"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.73%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "67.465.32"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.97%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "437.16"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.20%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.39"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.44%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "14.49"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.15%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "87.61"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.59"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.81%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "38.81"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.91%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.32"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.91%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.77"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.57%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "720.95"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.133"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "13.53"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.82"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.68%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "42.64"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.19%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "57.92"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("E35").Value = "  -0.71%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0₃0804"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +18.99%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("E41").Value = "  +6.74%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -7.97%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.82"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.49%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "147.62"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.80%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("E50").Value = "  +1.80%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "25.95"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.90%  "
